$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The practice data table (dates / fruit names / quantities) moves from
# A1:C7 down-and-right to D4:F10 - same values, same per-column styles.
$ws.Range("A1:C7").Cut($ws.Range("D4"))

# The source range is now empty; drop its leftover formatting too so the
# old A1:C7 cells go back to being plain/default cells.
$ws.Range("A1:C7").Clear()

# The custom column width that used to live on column A (the date column)
# now belongs to column D, the new home of the date column.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Selection moves along with the data.
$ws.Range("C8").Select() | Out-Null
